$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.833.23'
$ws.Range('E2').Value = '  -0.83%  '
$ws.Range('D3').Value = '2.034.32'
$ws.Range('E3').Value = '  -1.18%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.16'
$ws.Range('E5').Value = '  -1.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.614'
$ws.Range('E6').Value = '  -0.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.31'
$ws.Range('E7').Value = '  +3.76%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.385'
$ws.Range('E9').Value = '  -0.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0814'
$ws.Range('E10').Value = '  +1.05%  '
$ws.Range('E11').Value = '  +0.09%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.61'
$ws.Range('E12').Value = '  -0.02%  '
$ws.Range('D13').Value = '2.330.95'
$ws.Range('E13').Value = '  -1.42%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.05'
$ws.Range('E14').Value = '  +2.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.760'
$ws.Range('E15').Value = '  +0.74%  '
$ws.Range('E16').Value = '  -1.65%  '
$ws.Range('D17').Value = '2.032.61'
$ws.Range('E17').Value = '  -1.37%  '
$ws.Range('D18').Value = '37.787.25'
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.05'
$ws.Range('E19').Value = '  -1.60%  '
$ws.Range('E20').Value = '  -0.21%  '
$ws.Range('D21').Value = '0.0₃0825'
$ws.Range('E21').Value = '  -0.79%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '225.06'
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('E24').Value = '  -2.53%  '
$ws.Range('E25').Value = '  -1.76%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.26'
$ws.Range('E26').Value = '  -0.64%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.57'
$ws.Range('E27').Value = '  -0.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.129'
$ws.Range('E28').Value = '  -3.60%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.96'
$ws.Range('E29').Value = '  -0.53%  '
$ws.Range('E30').Value = '  -5.74%  '
$ws.Range('E31').Value = '  +1.32%  '
$ws.Range('E32').Value = '  -2.66%  '
$ws.Range('E33').Value = '  +4.17%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.51'
$ws.Range('E34').Value = '  -2.51%  '
$ws.Range('E35').Value = '  -2.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.36'
$ws.Range('E36').Value = '  +4.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.25'
$ws.Range('E37').Value = '  -4.33%  '
$ws.Range('E38').Value = '  -2.02%  '
$ws.Range('E39').Value = '  -0.14%  '
$ws.Range('D40').Value = '1.538.61'
$ws.Range('E40').Value = '  +3.94%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0218'
$ws.Range('E41').Value = '  -0.81%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '97.07'
$ws.Range('E42').Value = '  -1.38%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '16.89'
$ws.Range('E43').Value = '  -0.57%  '
$ws.Range('E44').Value = '  -0.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0922'
$ws.Range('E45').Value = '  -2.31%  '
$ws.Range('E46').Value = '  -1.43%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.93'
$ws.Range('E47').Value = '  -4.70%  '
$ws.Range('E48').Value = '  -1.04%  '
$ws.Range('E49').Value = '  +0.77%  '
$ws.Range('E50').Value = '  -0.27%  '
$ws.Range('D51').Value = '2.225.23'
$ws.Range('E51').Value = '  -1.22%  '
